$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.972.75"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "1.593.25"
$ws.Range("E3").Value = "  +0.17%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'210.41"
$ws.Range("E5").Value = "  +0.21%  "

$ws.Range("E7").Value = "  -0.23%  "

$ws.Range("E8").Value = "  -0.99%  "

$ws.Range("E9").Value = "  -1.21%  "

$ws.Range("D10").Value = "'17.94"
$ws.Range("E10").Value = "  -1.44%  "

$ws.Range("D11").Value = "'0.0808"
$ws.Range("E11").Value = "  +2.41%  "

$ws.Range("D12").Value = "1.816.77"
$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("D13").Value = "1.596.89"
$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").Value = "'3.99"
$ws.Range("E14").Value = "  -1.16%  "

$ws.Range("E15").Value = "  +0.09%  "

$ws.Range("D16").Value = "25.976.35"
$ws.Range("E16").Value = "  +0.34%  "

$ws.Range("D17").Value = "'60.06"
$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("D18").Value = "0.0₃0721"
$ws.Range("E18").Value = "  -0.21%  "

$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").Value = "'199.93"
$ws.Range("E20").Value = "  +3.71%  "

$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").Value = "'9.23"
$ws.Range("E22").Value = "  -1.89%  "

$ws.Range("E23").Value = "  +0.98%  "

$ws.Range("D24").Value = "'1.87"
$ws.Range("E24").Value = "  +9.82%  "

$ws.Range("D25").Value = "'143.12"
$ws.Range("E25").Value = "  +1.19%  "

$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("E27").Value = "  -8.21%  "

$ws.Range("D28").Value = "'15.07"
$ws.Range("E28").Value = "  -0.51%  "

$ws.Range("E29").Value = "  -0.45%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("E31").Value = "  +0.46%  "

$ws.Range("D32").Value = "'3.12"
$ws.Range("E32").Value = "  -0.14%  "

$ws.Range("D33").Value = "'2.95"
$ws.Range("E33").Value = "  -3.00%  "

$ws.Range("E34").Value = "  -1.52%  "

$ws.Range("D35").Value = "'2.37"
$ws.Range("E35").Value = "  +0.43%  "

$ws.Range("D36").Value = "1.121.97"
$ws.Range("E36").Value = "  +1.31%  "

$ws.Range("D37").Value = "'0.0163"
$ws.Range("E37").Value = "  +8.61%  "

$ws.Range("E38").Value = "  -0.15%  "

$ws.Range("E39").Value = "  -1.14%  "

$ws.Range("D40").Value = "'0.781"
$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("D41").Value = "'0.491"
$ws.Range("E41").Value = "  -2.90%  "

$ws.Range("D42").Value = "'0.785"
$ws.Range("E42").Value = "  -4.15%  "

$ws.Range("D43").Value = "1.727.88"
$ws.Range("E43").Value = "  +0.19%  "

$ws.Range("D44").Value = "'5.10"
$ws.Range("E44").Value = "  -1.58%  "

$ws.Range("D45").Value = "'92.55"
$ws.Range("E45").Value = "  -1.27%  "

$ws.Range("E46").Value = "  -1.34%  "

$ws.Range("D47").Value = "'53.29"
$ws.Range("E47").Value = "  -0.39%  "

$ws.Range("E48").Value = "  -1.30%  "

$ws.Range("E49").Value = "  -0.20%  "

$ws.Range("E50").Value = "  +0.22%  "

$ws.Range("D51").Value = "0.0₇0918"
$ws.Range("E51").Value = "  -17.26%  "
